$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 632.5294
$ws.Range("I39").Value = 410.84616
$ws.Range("J39").Value = 1353
$ws.Range("K39").Value = 1232.53848
$ws.Range("L39").Value = 4059
$ws.Range("M39").Value = -936.5384799999999
$ws.Range("N39").Value = -4651
$ws.Range("H69").Value = 3666.3333
$ws.Range("I69").Value = 2999.5
$ws.Range("K69").Value = 8998.5
$ws.Range("M69").Value = -8124.5
$ws.Range("H72").Value = 3666.3333
$ws.Range("I72").Value = 2999.5
$ws.Range("K72").Value = 26995.5
$ws.Range("M72").Value = -22627.5
$ws.Range("H96").Value = 554.5357
$ws.Range("I96").Value = 406.9565
$ws.Range("J96").Value = 1233.4
$ws.Range("K96").Value = 1220.8695
$ws.Range("L96").Value = 3700.2
$ws.Range("M96").Value = 152.1305
$ws.Range("N96").Value = -6446.200000000001
$ws.Range("H98").Value = 1832
$ws.Range("I98").Value = 1874.2727
$ws.Range("K98").Value = 1874.2727
$ws.Range("M98").Value = -376.2727
$ws.Range("H103").Value = 1836.6
$ws.Range("I103").Value = 1858.25
$ws.Range("K103").Value = 5574.75
$ws.Range("M103").Value = -4988.75
$ws.Range("H111").Value = 1589
$ws.Range("I111").Value = 1589
$ws.Range("K111").Value = 4767
$ws.Range("M111").Value = -1700
$ws.Range("H122").Value = 1832
$ws.Range("I122").Value = 1874.2727
$ws.Range("K122").Value = 5622.8181
$ws.Range("M122").Value = -3172.8181
$ws.Range("H141").Value = 3091.1667
$ws.Range("I141").Value = 2959.5
$ws.Range("K141").Value = 8878.5
$ws.Range("M141").Value = -3698.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4078.027
$ws.Range("I2").Value = 3253.2334
$ws.Range("J2").Value = 7612.857
$ws.Range("K2").Value = 3253.2334
$ws.Range("L2").Value = 7612.857
$ws.Range("M2").Value = -3140.2334
$ws.Range("N2").Value = -7838.857
$ws.Range("H45").Value = 9384.3125
$ws.Range("I45").Value = 10197.77
$ws.Range("K45").Value = 10197.77
$ws.Range("M45").Value = -9820.77
$ws.Range("H74").Value = 2197.26
$ws.Range("I74").Value = 1565.6857
$ws.Range("K74").Value = 1565.6857
$ws.Range("M74").Value = -691.6857
$ws.Range("H77").Value = 2197.26
$ws.Range("I77").Value = 1565.6857
$ws.Range("K77").Value = 7828.4285
$ws.Range("M77").Value = -3460.4285
$ws.Range("H88").Value = 2883.8333
$ws.Range("J88").Value = 3079.4
$ws.Range("L88").Value = 3079.4
$ws.Range("N88").Value = -3891.4
$ws.Range("H91").Value = 2883.8333
$ws.Range("J91").Value = 3079.4
$ws.Range("L91").Value = 3079.4
$ws.Range("N91").Value = -5887.4
$ws.Range("H110").Value = 846.9
$ws.Range("I110").Value = 722.2857
$ws.Range("K110").Value = 722.2857
$ws.Range("M110").Value = 1322.7143
$ws.Range("H116").Value = 4078.027
$ws.Range("I116").Value = 3253.2334
$ws.Range("J116").Value = 7612.857
$ws.Range("K116").Value = 3253.2334
$ws.Range("L116").Value = 7612.857
$ws.Range("M116").Value = -959.2334000000001
$ws.Range("N116").Value = -12200.857

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4078.027
$ws.Range("I3").Value = 3253.2334
$ws.Range("J3").Value = 7612.857
$ws.Range("K3").Value = 3253.2334
$ws.Range("L3").Value = 7612.857
$ws.Range("M3").Value = -3139.2334
$ws.Range("N3").Value = -7840.857
$ws.Range("H44").Value = 34883.223
$ws.Range("J44").Value = 34883.223
$ws.Range("L44").Value = 34883.223
$ws.Range("N44").Value = -35877.223
$ws.Range("H107").Value = 1271.3077
$ws.Range("I107").Value = 1143.0454
$ws.Range("J107").Value = 1976.75
$ws.Range("K107").Value = 1143.0454
$ws.Range("L107").Value = 1976.75
$ws.Range("M107").Value = 776.9546
$ws.Range("N107").Value = -5816.75
$ws.Range("H134").Value = 7000.946
$ws.Range("I134").Value = 3919.926
$ws.Range("J134").Value = 15319.7
$ws.Range("K134").Value = 11759.778
$ws.Range("L134").Value = 45959.10000000001
$ws.Range("M134").Value = -9224.778
$ws.Range("N134").Value = -51029.10000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3131
$ws.Range("I31").Value = 1363
$ws.Range("J31").Value = 5959.8
$ws.Range("K31").Value = 1363
$ws.Range("L31").Value = 5959.8
$ws.Range("M31").Value = -1068
$ws.Range("N31").Value = -6549.8
$ws.Range("H34").Value = 3131
$ws.Range("I34").Value = 1363
$ws.Range("J34").Value = 5959.8
$ws.Range("K34").Value = 1363
$ws.Range("L34").Value = 5959.8
$ws.Range("M34").Value = -1161
$ws.Range("N34").Value = -6363.8
$ws.Range("H134").Value = 9018.225
$ws.Range("I134").Value = 7107.7095
$ws.Range("J134").Value = 15598.889
$ws.Range("K134").Value = 21323.1285
$ws.Range("L134").Value = 46796.667
$ws.Range("M134").Value = -18788.1285
$ws.Range("N134").Value = -51866.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 861.8889
$ws.Range("I113").Value = 1430.6666
$ws.Range("J113").Value = 577.5
$ws.Range("K113").Value = 4291.9998
$ws.Range("L113").Value = 1732.5
$ws.Range("M113").Value = -2121.9998
$ws.Range("N113").Value = -6072.5
$ws.Range("H123").Value = 1492.5
$ws.Range("I123").Value = 1492.5
$ws.Range("K123").Value = 4477.5
$ws.Range("M123").Value = -2027.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7101.1665
$ws.Range("I70").Value = 6303.75
$ws.Range("J70").Value = 7898.5835
$ws.Range("K70").Value = 6303.75
$ws.Range("L70").Value = 7898.5835
$ws.Range("M70").Value = -6033.75
$ws.Range("N70").Value = -8438.583500000001
$ws.Range("H73").Value = 7101.1665
$ws.Range("I73").Value = 6303.75
$ws.Range("J73").Value = 7898.5835
$ws.Range("K73").Value = 6303.75
$ws.Range("L73").Value = 7898.5835
$ws.Range("M73").Value = -5367.75
$ws.Range("N73").Value = -9770.583500000001
$ws.Range("H113").Value = 225428.67
$ws.Range("I113").Value = 335059.66
$ws.Range("J113").Value = 6166.6665
$ws.Range("K113").Value = 335059.66
$ws.Range("L113").Value = 6166.6665
$ws.Range("M113").Value = -332889.66
$ws.Range("N113").Value = -10506.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1027.8096
$ws.Range("I16").Value = 1027.8096
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1027.8096
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -857.8096
$ws.Range("N16").ClearContents()
$ws.Range("H46").Value = 1214.3103
$ws.Range("I46").Value = 990.375
$ws.Range("J46").Value = 1299.619
$ws.Range("K46").Value = 990.375
$ws.Range("L46").Value = 1299.619
$ws.Range("M46").Value = -802.375
$ws.Range("N46").Value = -1675.619
$ws.Range("H61").Value = 1338.4286
$ws.Range("I61").Value = 1182.4375
$ws.Range("K61").Value = 1182.4375
$ws.Range("M61").Value = -980.4375
$ws.Range("H113").Value = 1338.4286
$ws.Range("I113").Value = 1182.4375
$ws.Range("K113").Value = 1182.4375
$ws.Range("M113").Value = 987.5625
$ws.Range("H136").Value = 2695.139
$ws.Range("I136").Value = 2413.3
$ws.Range("J136").Value = 4104.3335
$ws.Range("K136").Value = 7239.900000000001
$ws.Range("L136").Value = 12313.0005
$ws.Range("M136").Value = -4689.900000000001
$ws.Range("N136").Value = -17413.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 3498.5
$ws.Range("I29").Value = 997
$ws.Range("J29").Value = 6000
$ws.Range("K29").Value = 997
$ws.Range("L29").Value = 6000
$ws.Range("M29").Value = -707
$ws.Range("N29").Value = -6580
$ws.Range("H42").Value = 66665.664
$ws.Range("I42").Value = 49999
$ws.Range("K42").Value = 49999
$ws.Range("M42").Value = -49621
$ws.Range("H107").Value = 1616.0358
$ws.Range("I107").Value = 1853.6666
$ws.Range("J107").Value = 1437.8125
$ws.Range("K107").Value = 5560.9998
$ws.Range("L107").Value = 4313.4375
$ws.Range("M107").Value = -3640.9998
$ws.Range("N107").Value = -8153.4375
$ws.Range("H132").Value = 20718.52
$ws.Range("I132").Value = 14838.467
$ws.Range("K132").Value = 44515.401
$ws.Range("M132").Value = -41985.401
